# "added diversity to the plot." — the author stripped the per-species
# "Notes (areas/fire response)" (col H) and "Conservation status" (col I)
# annotations out of the Relevant_Species summary table so the sheet no
# longer hard-codes a single status per species (those notes were no
# longer representative once more diversity/abundance data was folded
# into the plot). Clear the values in H2:I43 but keep the column
# formatting/styles intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relevant_Species")

$ws.Range("H2:I43").ClearContents()
